$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Headers for the new "anomaly" / "lat" helper table
$ws.Range("N1").Value = "anomaly"
$ws.Range("O1").Value = "lat"

# Row 2 (single-row formulas, not shared)
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 90
$ws.Range("P2").Formula = "=N2-180"
$ws.Range("Q2").Formula = "=MOD(P2,90)"
$ws.Range("R2").Formula = "=IF(N2>180,-MOD(N2,180),N2)"
$ws.Range("T2").Formula = "=90-N2"

# Row 3 values
$ws.Range("N3").Value = 45
$ws.Range("O3").Value = 45

# Row 4 values (no N/O values given explicitly besides formulas row - N4/O4 present)
$ws.Range("N4").Value = 90
$ws.Range("O4").Value = 0

# Row 5 values
$ws.Range("N5").Value = 120
$ws.Range("O5").Value = -30

# Row 6 (brand new row) values
$ws.Range("N6").Value = 180
$ws.Range("O6").Value = -90

# Row 7 values
$ws.Range("N7").Value = 220
$ws.Range("O7").Value = -50

# Row 8 values
$ws.Range("N8").Value = 270
$ws.Range("O8").Value = 0

# Row 9 values
$ws.Range("N9").Value = 330
$ws.Range("O9").Value = 60

# Shared formulas: P3:P9, Q3:Q9, R3:R9, T3:T6, T8:T9
$ws.Range("P3:P9").Formula = "=N3-180"
$ws.Range("Q3:Q9").Formula = "=MOD(P3,90)"
$ws.Range("R3:R9").Formula = "=IF(N3>180,-MOD(N3,180),N3)"
$ws.Range("T3:T6").Formula = "=90-N3"

# Row 7's T cell keeps its own (non-shared) formula matching the T8:T9 pattern
$ws.Range("T7").Formula = "=-90+MOD(N7,180)"
$ws.Range("T8:T9").Formula = "=-90+MOD(N8,180)"

# Rows 13-25: Q13/R13 become the masters of new shared-formula groups spanning Q13:Q25 / R13:R25
$ws.Range("Q13:Q25").Formula = "=MOD(P13+90, 360)"
$ws.Range("R13:R25").Formula = "=MOD(P13+270, 360)"

# Update the active cell / selection
$ws.Range("V6").Select() | Out-Null
